# Apply updates to the COVID closures by state workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Arkansas (row 5): add "complete" date (column F)
$ws.Range("F5").Value = 43925
$ws.Range("D5").Copy()
$ws.Range("F5").PasteSpecial(-4122)

# Mississippi (row 26): add "complete" date (column F)
$ws.Range("F26").Value = 43914
$ws.Range("D26").Copy()
$ws.Range("F26").PasteSpecial(-4122)

# Missouri (row 27): add "complete" date (column F)
$ws.Range("F27").Value = 43927
$ws.Range("D27").Copy()
$ws.Range("F27").PasteSpecial(-4122)

# North Dakota (row 36): fix "partial" date, add blank styled "complete" cell
$ws.Range("E36").Value = 43909
$ws.Range("D36").Copy()
$ws.Range("F36").PasteSpecial(-4122)

# Update the sheet view: clear frozen/top-left scroll position, move selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E37").Select() | Out-Null
